# Update the scraped "想去人数" (want-to-go count) and one "最低票价"
# (min ticket price) figures across the workbook's sheets, matching the
# refreshed gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 8091
$ws1.Range("F4").Value  = 1904
$ws1.Range("F5").Value  = 6492
$ws1.Range("F6").Value  = 157
$ws1.Range("F7").Value  = 2042
$ws1.Range("F8").Value  = 559
$ws1.Range("F9").Value  = 40
$ws1.Range("F14").Value = 64
$ws1.Range("F15").Value = 8435
$ws1.Range("G15").Value = 75
$ws1.Range("F17").Value = 62
$ws1.Range("F19").Value = 113
$ws1.Range("F20").Value = 1799
$ws1.Range("F25").Value = 17
$ws1.Range("F28").Value = 2
$ws1.Range("F30").Value = 2030
$ws1.Range("F31").Value = 841
$ws1.Range("F32").Value = 461
$ws1.Range("F35").Value = 164
$ws1.Range("F36").Value = 141
$ws1.Range("F37").Value = 4

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 706

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 706
$ws4.Range("F6").Value  = 8091
$ws4.Range("F9").Value  = 1904
$ws4.Range("F10").Value = 6492
$ws4.Range("F11").Value = 2043
$ws4.Range("F13").Value = 559
$ws4.Range("F14").Value = 40
$ws4.Range("F22").Value = 64
$ws4.Range("F23").Value = 8435
$ws4.Range("G23").Value = 75
$ws4.Range("F25").Value = 62
$ws4.Range("F27").Value = 113
$ws4.Range("F28").Value = 1799
$ws4.Range("F35").Value = 2030
$ws4.Range("F36").Value = 841
$ws4.Range("F38").Value = 461
$ws4.Range("F41").Value = 141
